$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 6 (RCB vs KKR) - row 18 - enter the scores for each player
$ws.Range("E18").Value = 50
$ws.Range("H18").Value = 70
$ws.Range("K18").Value = 80
$ws.Range("N18").Value = 40
$ws.Range("Q18").Value = 100
$ws.Range("T18").Value = 60
$ws.Range("W18").Value = 0
